$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.631.51"
$ws.Range("E2").Value = "  -5.18%  "
$ws.Range("D3").Value = "2.582.38"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.575"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.14%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.552"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.56%  "
$ws.Range("D13").Value = "2.968.21"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.107"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("D15").Value = "2.584.36"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.883"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.71%  "
$ws.Range("D18").Value = "43.518.40"
$ws.Range("E18").Value = "  -5.74%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0971"
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "264.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.66%  "
$ws.Range("E29").Value = "  -4.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("E34").Value = "  -1.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "151.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0806"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.117"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.120"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0313"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.85%  "
$ws.Range("D44").Value = "2.071.57"
$ws.Range("E44").Value = "  -3.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.996"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.73%  "
$ws.Range("E48").Value = "  +3.60%  "
$ws.Range("D49").Value = "2.836.79"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.190"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.61%  "
